# Remove the stray "∆" (U+2206 INCREMENT) run from the "Phase deviation
# ∆φsrc" / "Phase deviation ∆φtrg" text boxes ("Textfeld 64") on the
# slide. Everything else in those runs (the "φ" and the "src"/"trg"
# suffix, each carrying their own run-level formatting) is left intact.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$delta = [char]0x2206   # "∆"

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)

    if (-not $sh.HasTextFrame) { continue }
    if ($sh.Name -ne "Textfeld 64") { continue }

    $tr = $sh.TextFrame.TextRange
    $pos = $tr.Text.IndexOf($delta)
    if ($pos -lt 0) { continue }

    # Record the shape's current (pre-edit) height so the autosize
    # textbox keeps its original extent after the text mutation -
    # PowerPoint re-lays out spAutoFit boxes on any text edit.
    $origHeight = $sh.Height

    # IndexOf is 0-based; Characters() is 1-based -> +1.
    $run = $tr.Characters($pos + 1, 1)
    $run.Text = ""

    $sh.Height = $origHeight + 0.00001
}
